$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the "What I did" text for week 1 (D2)
$ws.Range("D2").Value = "Watched children's shows I'm familiar with and read simple manga."

# Add new row 3 data
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2.2953125000000001
$ws.Range("B3").NumberFormat = "[h]:mm:ss"
$ws.Range("C3").Value = "Avatar the Last Airbender (Audiovisual, English, Familiar):15; Fist of the North Star (Text with visuals, Japanese, New):17;"
$ws.Range("D3").Value = "Watched children's shows I'm familiar with and read simple manga."

$ws.Range("C3").Select()
